# "Generate Report for Handoff" - refresh the localization-status report:
#   - the ea844430-... file (row 3 on every sheet) is now ready to be handed
#     off, so its Status flips from "Handed back: in sync with en-US" to
#     "Ready for handoff"
#   - the Latest Handoff Datetime for that same file moved forward
#     (the zh-cn / de-de target languages were handed off again)

$wb = $excel.ActiveWorkbook

$statusOld = "Handed back: in sync with en-US"
$statusNew = "Ready for handoff"

# ---- Overview sheet ---------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = $statusNew
$overview.Range("C3").Value = $statusNew
# Latest Handoff Date for the ea844430 entry (shared with row 2, same value)
$overview.Range("D2").Value = "2016-03-24 03:12:40"
$overview.Range("D3").Value = "2016-03-24 03:12:40"

# ---- zh-cn sheet --------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = $statusNew
$zhcn.Range("E2").Value = "2016-03-24 03:12:32"
$zhcn.Range("E3").Value = "2016-03-24 03:12:32"

# ---- de-de sheet --------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = $statusNew
$dede.Range("E2").Value = "2016-03-24 03:12:40"
$dede.Range("E3").Value = "2016-03-24 03:12:40"
